$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newValues = @(
    5295.952345967664,
    -2147.142791632776,
    9923.437984105742,
    12553.89319920146,
    5048.257990149968,
    8176.12212012359,
    16336.68841752559,
    4047.442791702304,
    1331.281148206443,
    4814.798501938374,
    22195.84357116088,
    7413.673506962974,
    10716.53772514699,
    16349.17501649283,
    5911.063232745426,
    14477.67432737427,
    6643.079073157605,
    17119.6417362342,
    14269.6677016509,
    17042.30426069828,
    -264.2898158340777,
    24036.41859544036,
    13171.53619885951,
    4558.934979593558,
    11630.15462398761,
    11861.13432540124,
    23972.05851730144,
    9989.716149974167,
    9766.57921435846,
    10714.95876990029,
    18509.44742569692,
    7730.262519602175,
    9545.214252369886,
    9418.703531776073,
    18816.58785067385,
    9865.04896910398,
    19668.79068072772,
    8717.679848784412,
    21360.07318259204,
    15931.89724942467,
    630.8470349329073,
    9377.982062821025,
    15538.02980937499,
    20493.77539523179,
    10755.44781824835,
    12532.92569838882,
    3267.604639950074,
    25509.29035583191,
    26735.80215244205,
    6621.859719786033,
    1585.133438806324,
    7392.002497859261,
    19684.29486243112,
    24046.86781587759,
    24507.41777166135,
    18846.23187053042,
    14803.31722238317,
    4867.971897739782,
    1009.777906319037,
    22290.95883401286
)

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $newValues[$i]
}

